# Adjusting one table mapping
#
# Table 13 ("tab:ratios") was mapped to the wrong LaTeX source file.
# Row 27 held "text/analysis/table_ratios_edited.tex" (table number 13)
# while row 28 held "text/analysis/table_ratios.tex" with no table number.
# The table number "13" actually belongs on the table_ratios.tex row, so
# move it down one row (A27 -> A28), leaving the file/label columns as
# they were.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the table number out of row 27 ...
$ws.Range("A27").ClearContents()
# ... and onto row 28, where it belongs.
$ws.Range("A28").Value = 13

# Row 27 renders slightly taller once it no longer carries a number.
$ws.Rows("27").RowHeight = 12.8

# Leave the cursor where the edit was made.
$null = $ws.Range("F28").Select()
